$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The "Website / Source" column (D) used to show a short domain-style label
# (e.g. "Durable.co", "TeleportHQ.io", "Wix.com", ...) for each tool, while
# the hyperlink attached to the cell already pointed at the real URL. The
# sheet is being updated so the visible cell text is the actual URL itself
# (matching the URL shown for the tool immediately above it in a couple of
# cases where two rows share the same source link). The hyperlinks
# themselves keep pointing at the same addresses as before.
# ---------------------------------------------------------------------------

$rowUrls = @{
    3  = "https://durable.co/"
    4  = "https://durable.co/"
    5  = "https://teleporthq.io/"
    6  = "https://teleporthq.io/"
    7  = "https://www.wix.com/ai-website-builder"
    8  = "https://www.wix.com/ai-website-builder"
    9  = "https://www.appypie.com/"
    10 = "https://thunkable.com/"
    11 = "https://www.adalo.com/"
    12 = "https://www.adalo.com/"
    13 = "https://www.v.one/"
    14 = "https://www.v.one/"
    15 = "https://uizard.io/"
    16 = "https://bubble.io/"
    17 = "https://webflow.com/"
    18 = "https://www.figma.com/community"
    19 = "https://www.codux.io/"
    20 = "https://www.codux.io/"
    21 = "https://www.pineapplebuilder.com/"
    22 = "https://www.pineapplebuilder.com/"
}

for ($row = 3; $row -le 22; $row++) {
    $ws.Range("D$row").Value = $rowUrls[$row]
}
